$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price (D) column cells to remain text, since their values
# look numeric (e.g. '43.935.45', '0.678') and Excel would otherwise
# silently convert them to real numbers on assignment.
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = '@'
    $rng.Value = $value
    $rng.Style = 'Normal'
}

Set-TextValue 'D2' '43.935.45'
$ws.Range('E2').Value = '  -0.52%  '
Set-TextValue 'D3' '2.351.68'
$ws.Range('E3').Value = '  -0.30%  '
$ws.Range('E4').Value = '  +0.02%  '
Set-TextValue 'D5' '0.678'
$ws.Range('E5').Value = '  +0.05%  '
Set-TextValue 'D6' '238.72'
$ws.Range('E6').Value = '  +0.06%  '
Set-TextValue 'D7' '73.71'
$ws.Range('E7').Value = '  +0.56%  '
$ws.Range('E8').Value = '  -0.07%  '
Set-TextValue 'D9' '0.596'
$ws.Range('E9').Value = '  +9.01%  '
Set-TextValue 'D10' '0.100'
$ws.Range('E10').Value = '  -0.46%  '
Set-TextValue 'D11' '57.27'
$ws.Range('E11').Value = '  +0.13%  '
Set-TextValue 'D12' '32.34'
$ws.Range('E12').Value = '  +9.62%  '
Set-TextValue 'D13' '7.30'
$ws.Range('E13').Value = '  +8.64%  '
$ws.Range('E14').Value = '  +0.20%  '
Set-TextValue 'D15' '2.700.07'
$ws.Range('E15').Value = '  -0.51%  '
Set-TextValue 'D16' '16.56'
$ws.Range('E16').Value = '  -1.81%  '
Set-TextValue 'D17' '0.896'
$ws.Range('E17').Value = '  -0.88%  '
Set-TextValue 'D18' '2.350.00'
$ws.Range('E18').Value = '  -0.54%  '
Set-TextValue 'D19' '43.852.77'
$ws.Range('E19').Value = '  -0.58%  '
$ws.Range('E20').Value = '  -1.13%  '
Set-TextValue 'D21' '6.71'
$ws.Range('E21').Value = '  +3.93%  '
Set-TextValue 'D22' '76.74'
Set-TextValue 'D23' '257.76'
$ws.Range('E23').Value = '  +0.75%  '
$ws.Range('E24').Value = '  +21.31%  '
$ws.Range('E25').Value = '  -0.01%  '
$ws.Range('E27').Value = '  -1.51%  '
Set-TextValue 'D28' '10.72'
Set-TextValue 'D29' '2.28'
$ws.Range('E29').Value = '  -0.45%  '
Set-TextValue 'D30' '22.61'
$ws.Range('E30').Value = '  +0.67%  '
Set-TextValue 'D31' '175.39'
$ws.Range('E31').Value = '  +1.41%  '
Set-TextValue 'D32' '0.129'
$ws.Range('E32').Value = '  -2.81%  '
$ws.Range('E33').Value = '  +2.70%  '
Set-TextValue 'D34' '0.0763'
$ws.Range('E34').Value = '  +3.72%  '
Set-TextValue 'D35' '5.20'
$ws.Range('E35').Value = '  -0.21%  '
Set-TextValue 'D36' '5.44'
$ws.Range('E36').Value = '  +3.87%  '
Set-TextValue 'D37' '3.75'
$ws.Range('E37').Value = '  -4.59%  '
Set-TextValue 'D38' '2.35'
$ws.Range('E38').Value = '  -3.84%  '
Set-TextValue 'D39' '6.29'
$ws.Range('E39').Value = '  -2.94%  '
$ws.Range('E40').Value = '  +2.30%  '
Set-TextValue 'D41' '0.110'
$ws.Range('E41').Value = '  +11.64%  '
$ws.Range('E42').Value = '  +12.65%  '
Set-TextValue 'D43' '9.01'
$ws.Range('E43').Value = '  +1.79%  '
Set-TextValue 'D44' '18.91'
$ws.Range('E44').Value = '  -4.00%  '
$ws.Range('E45').Value = '  +0.01%  '
Set-TextValue 'D46' '4.73'
$ws.Range('E46').Value = '  +5.10%  '
$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue 'D47' '2.51'
$ws.Range('E47').Value = '  +6.26%  '
$ws.Range('B48').Value = 'MultiversX'
$ws.Range('C48').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
Set-TextValue 'D48' '57.89'
$ws.Range('E48').Value = '  +9.47%  '
Set-TextValue 'D49' '1.24'
$ws.Range('E49').Value = '  -0.67%  '
$ws.Range('E50').Value = '  -0.06%  '
Set-TextValue 'D51' '99.83'
$ws.Range('E51').Value = '  +1.15%  '
